# Daily attendance processing - 2026-01-20 15:16:11
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 194: session became Recorded (copy "Recorded" row formatting from row 193, then set new values)
$ws.Range("A193:I193").Copy()
$ws.Range("A194:I194").PasteSpecial(-4122)
$ws.Range("G194").Value = "nadahassanein99@hotmail.com"
$ws.Range("H194").Value = "40/63"
$ws.Range("I194").Value = "Recorded"

# Remaining cell value updates
$ws.Range("H2").Value = "60/67"
$ws.Range("H3").Value = "0/67"
$ws.Range("H4").Value = "66/67"
$ws.Range("H5").Value = "51/67"
$ws.Range("H6").Value = "0/67"
$ws.Range("L6").Value = 204
$ws.Range("H7").Value = "59/67"
$ws.Range("L7").Value = 33
$ws.Range("H8").Value = "60/67"
$ws.Range("H9").Value = "66/67"
$ws.Range("L9").Value = "74.2%"
$ws.Range("H10").Value = "40/67"
$ws.Range("L10").Value = "74.9%"
$ws.Range("H11").Value = "27/67"
$ws.Range("H12").Value = "53/67"
$ws.Range("H13").Value = "0/67"
$ws.Range("H14").Value = "47/67"
$ws.Range("H15").Value = "25/67"
$ws.Range("M15").Value = 67
$ws.Range("S15").Value = "71.8%"
$ws.Range("H16").Value = "51/67"
$ws.Range("H17").Value = "58/67"
$ws.Range("H18").Value = "58/67"
$ws.Range("O18").Value = 43
$ws.Range("P18").Value = 6
$ws.Range("R18").Value = "78.2%"
$ws.Range("S18").Value = "76.5%"
$ws.Range("H19").Value = "0/67"
$ws.Range("M19").Value = 61
$ws.Range("S19").Value = "69.2%"
$ws.Range("H20").Value = "28/67"
$ws.Range("H21").Value = "31/67"
$ws.Range("H22").Value = "45/67"
$ws.Range("H23").Value = "38/67"
$ws.Range("H24").Value = "45/67"
$ws.Range("H25").Value = "0/67"
$ws.Range("G26").Value = "haderreda2919@gmail.com, emp17.nada.h.attia@gmail.com"
$ws.Range("H26").Value = "66/67"
$ws.Range("H27").Value = "62/67"
$ws.Range("H28").Value = "55/67"
$ws.Range("H29").Value = "59/67"
$ws.Range("H30").Value = "0/67"
$ws.Range("H31").Value = "0/67"
$ws.Range("H32").Value = "57/67"
$ws.Range("G33").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("H33").Value = "60/67"
$ws.Range("H34").Value = "52/67"
$ws.Range("H35").Value = "30/67"
$ws.Range("H36").Value = "0/67"
$ws.Range("G37").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("H37").Value = "45/67"
$ws.Range("H38").Value = "57/67"
$ws.Range("H39").Value = "59/67"
$ws.Range("H40").Value = "27/67"
$ws.Range("H41").Value = "17/67"
$ws.Range("H42").Value = "47/67"
$ws.Range("H43").Value = "54/67"
$ws.Range("H44").Value = "60/67"
$ws.Range("H45").Value = "55/67"
$ws.Range("H46").Value = "42/67"
$ws.Range("H47").Value = "1/67"
$ws.Range("H48").Value = "59/67"
$ws.Range("H49").Value = "52/67"
$ws.Range("H50").Value = "54/67"
$ws.Range("H51").Value = "43/67"
$ws.Range("H52").Value = "0/67"
$ws.Range("H53").Value = "0/67"
$ws.Range("H54").Value = "0/67"
$ws.Range("H55").Value = "0/67"
$ws.Range("H56").Value = "0/67"
$ws.Range("G89").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("G90").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("G95").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("G125").Value = "abdallahashraf2023@gmail.com, ahmedali78112@gmail.com"
$ws.Range("G197").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("G199").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("G202").Value = "160534@med.asu.edu.eg, emp17.mariam.m.goda@gmail.com"
$ws.Range("H222").Value = "43/61"
$ws.Range("H223").Value = "31/61"
$ws.Range("H224").Value = "50/61"
$ws.Range("H225").Value = "38/61"
$ws.Range("H226").Value = "0/61"
$ws.Range("H227").Value = "0/61"
$ws.Range("H228").Value = "55/61"
$ws.Range("H229").Value = "0/61"
$ws.Range("H230").Value = "50/61"
$ws.Range("H231").Value = "35/61"
$ws.Range("H232").Value = "40/61"
$ws.Range("H233").Value = "49/61"
$ws.Range("H234").Value = "54/61"
$ws.Range("H235").Value = "0/61"
$ws.Range("H236").Value = "0/61"
$ws.Range("H237").Value = "0/61"
$ws.Range("H238").Value = "0/61"
$ws.Range("H239").Value = "0/61"
$ws.Range("H240").Value = "0/61"
$ws.Range("H241").Value = "0/61"
$ws.Range("H242").Value = "37/61"
$ws.Range("H243").Value = "38/61"
$ws.Range("H244").Value = "31/61"
$ws.Range("H245").Value = "44/61"
$ws.Range("H246").Value = "0/61"
$ws.Range("H247").Value = "42/61"
$ws.Range("H248").Value = "51/61"
$ws.Range("H249").Value = "50/61"
$ws.Range("H250").Value = "43/61"
$ws.Range("H251").Value = "0/61"
$ws.Range("H252").Value = "22/61"
$ws.Range("H253").Value = "49/61"
$ws.Range("H254").Value = "29/61"
$ws.Range("H255").Value = "56/61"
$ws.Range("H256").Value = "23/61"
$ws.Range("H257").Value = "23/61"
$ws.Range("H258").Value = "51/61"
$ws.Range("H259").Value = "0/61"
$ws.Range("H260").Value = "43/61"
$ws.Range("H261").Value = "25/61"
$ws.Range("H262").Value = "35/61"
$ws.Range("H263").Value = "52/61"
$ws.Range("H264").Value = "42/61"
$ws.Range("H265").Value = "0/61"
$ws.Range("H266").Value = "54/61"
$ws.Range("H267").Value = "45/61"
$ws.Range("H268").Value = "58/61"
$ws.Range("H269").Value = "54/61"
$ws.Range("H270").Value = "22/61"
$ws.Range("H271").Value = "38/61"
$ws.Range("H272").Value = "48/61"
$ws.Range("H273").Value = "55/61"
$ws.Range("H274").Value = "42/61"
$ws.Range("H275").Value = "44/61"
$ws.Range("H276").Value = "40/61"
